{"js": "// Update the worksheet date header and the 25 two-digit / one-digit\n// division problems laid out in the 5x5 answer table (every 4th table\n// row holds a row of 5 answer cells; the rows between are spacers).\n\n// 1) Update the date/weekday heading paragraph.\nconst body = context.document.body;\nconst headingParas = body.paragraphs;\nheadingParas.load(\"items/text\");\nawait context.sync();\n\nconst oldHeading = \"2024-04-13 Saturday\";\nconst newHeading = \"2024-04-14 Sunday\";\nfor (const para of headingParas.items) {\n  if (para.text === oldHeading) {\n    para.insertText(newHeading, \"Replace\");\n    break;\n  }\n}\n\n// 2) Update the answer table's 25 cells (5 populated rows x 5 columns).\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// New answers, in the same row/column order as the populated table rows.\nconst newAnswers = [\n  [\"50\u00f73=16, 2\", \"20\u00f77=2, 6\", \"13\u00f72=6, 1\", \"68\u00f73=22, 2\", \"40\u00f78=5, 0\"],\n  [\"95\u00f75=19, 0\", \"90\u00f78=11, 2\", \"88\u00f77=12, 4\", \"37\u00f75=7, 2\", \"22\u00f73=7, 1\"],\n  [\"52\u00f79=5, 7\", \"12\u00f77=1, 5\", \"77\u00f79=8, 5\", \"30\u00f74=7, 2\", \"75\u00f73=25, 0\"],\n  [\"25\u00f74=6, 1\", \"22\u00f72=11, 0\", \"98\u00f77=14, 0\", \"99\u00f76=16, 3\", \"42\u00f73=14, 0\"],\n  [\"25\u00f75=5, 0\", \"99\u00f73=33, 0\", \"96\u00f75=19, 1\", \"58\u00f73=19, 1\", \"59\u00f79=6, 5\"],\n];\n\n// The answers live on every 4th row (0, 4, 8, 12, 16); the rows in\n// between are blank spacer rows (still made of 5 empty cells, so\n// cellCount can't tell them apart). Discover the populated rows by\n// checking the first cell's text instead of hard-coding the indices,\n// so the script is resilient to row-count differences.\ntable.load(\"rowCount\");\nawait context.sync();\n\nconst populatedRowIndices = [];\nfor (let i = 0; i < table.rowCount; i++) {\n  const firstCell = table.getCell(i, 0);\n  firstCell.body.load(\"text\");\n  await context.sync();\n  if (firstCell.body.text.trim().length > 0) {\n    populatedRowIndices.push(i);\n  }\n}\n\nfor (let r = 0; r < populatedRowIndices.length && r < newAnswers.length; r++) {\n  const rowIndex = populatedRowIndices[r];\n  const rowValues = newAnswers[r];\n  for (let c = 0; c < rowValues.length; c++) {\n    const cell = table.getCell(rowIndex, c);\n    cell.body.paragraphs.load(\"items\");\n    await context.sync();\n    cell.body.paragraphs.items[0].insertText(rowValues[c], \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the worksheet date header and the 25 two-digit / one-digit\n# division problems laid out in the 5x5 answer table (every 4th table\n# row holds a row of 5 answer cells; the rows between are spacers).\n\n$d = $word.ActiveDocument\n\n# 1) Update the date/weekday heading paragraph.\n$d.Content.Find.Execute(\n    \"2024-04-13 Saturday\", $false, $false, $false, $false, $false,\n    $true, 1, $false, \"2024-04-14 Sunday\", 2\n) | Out-Null\n\n# 2) Update the answer table's 25 cells (5 populated rows x 5 columns).\n$t = $d.Tables(1)\n\n# New answers, in the same row/column order as the populated table rows.\n$newAnswers = @(\n    @(\"50\u00f73=16, 2\", \"20\u00f77=2, 6\", \"13\u00f72=6, 1\", \"68\u00f73=22, 2\", \"40\u00f78=5, 0\"),\n    @(\"95\u00f75=19, 0\", \"90\u00f78=11, 2\", \"88\u00f77=12, 4\", \"37\u00f75=7, 2\", \"22\u00f73=7, 1\"),\n    @(\"52\u00f79=5, 7\", \"12\u00f77=1, 5\", \"77\u00f79=8, 5\", \"30\u00f74=7, 2\", \"75\u00f73=25, 0\"),\n    @(\"25\u00f74=6, 1\", \"22\u00f72=11, 0\", \"98\u00f77=14, 0\", \"99\u00f76=16, 3\", \"42\u00f73=14, 0\"),\n    @(\"25\u00f75=5, 0\", \"99\u00f73=33, 0\", \"96\u00f75=19, 1\", \"58\u00f73=19, 1\", \"59\u00f79=6, 5\")\n)\n\n# The answers live on every 4th row (1, 5, 9, 13, 17 - Word rows are\n# 1-indexed); the rows in between are blank spacer rows. Discover the\n# populated rows by checking the first cell's text instead of\n# hard-coding the indices, so the script is resilient to row-count\n# differences.\n$populatedRows = New-Object System.Collections.ArrayList\nfor ($i = 1; $i -le $t.Rows.Count; $i++) {\n    $cellText = $t.Cell($i, 1).Range.Text\n    $cellText = $cellText -replace \"[\\x07\\x0d\\x0c]\", \"\"\n    if ($cellText.Trim().Length -gt 0) {\n        [void]$populatedRows.Add($i)\n    }\n}\n\nfor ($r = 0; $r -lt $populatedRows.Count -and $r -lt $newAnswers.Count; $r++) {\n    $rowIndex = $populatedRows[$r]\n    $rowValues = $newAnswers[$r]\n    for ($c = 0; $c -lt $rowValues.Count; $c++) {\n        $t.Cell($rowIndex, $c + 1).Range.Text = $rowValues[$c]\n    }\n}\n"}
